$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "La Palma" and "Arroyo de la Luz" order (A56 <-> A57 display values)
$ws.Range("A56").Value = "La Palma"
$ws.Range("A57").Value = "Arroyo de la Luz"

# Update "Datos actualizados" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 21:46"

# Update Navarra (row 8) Recuperados/Muertes counts
$ws.Range("D8").Value = 582
$ws.Range("E8").Value = 9
